$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("counts")

$ws.Range("B4").Value = "crt:72_73_74:C/A_V/A_I"
$ws.Range("B5").Value = "crt:72_73_74:C|A|D_V|A|D_I"

$ws.Range("B5").Select() | Out-Null
